# Append 8 new roadway-link rows (216-223) to the end of the data table on
# Hoja1, matching the rows that the upstream FreightNetwork derivation
# added to data/roadway_links.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# id_link, distance, gauge
$newRows = @(
    @("17-1011",   199, "unica"),
    @("21-1021",   120, "unica"),
    @("1-1003",     58, "unica"),
    @("1037-1052",  56, "unica"),
    @("46-1037",    94, "unica"),
    @("49-1023",   170, "unica"),
    @("3-1003",     45, "unica"),
    @("1-1021",    185, "unica")
)

$firstNewRow = 216
$lastNewRow = $firstNewRow + $newRows.Count - 1

# Give column B (distance) the same look as the rest of the data column
# (font/number-format used by B2:B215) before the values go in, so the new
# cells pick up that formatting instead of the sheet's bare default style.
$ws.Range("B215").Copy() | Out-Null
$ws.Range("B" + $firstNewRow + ":B" + $lastNewRow).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$r = $firstNewRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

# Mirror the author's final cursor position: scrolled down near the bottom
# of the sheet with the last cell of the new block selected.
$ws.Range("C" + $lastNewRow).Select() | Out-Null
